# CORELIMS-98 - production addition of F3
# Update RCK rack barcodes in column C (LOCATION BARCODE) for freezer_03_shelf2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C8").Value = "RCK105"
$ws.Range("C9:C15").Value = "RCK106"
$ws.Range("C16:C22").Value = "RCK107"
$ws.Range("C23:C29").Value = "RCK108"
$ws.Range("C30:C36").Value = "RCK109"
